# Auto-generated Excel COM-interop script applying the Jenova_Profits numeric updates
# (scheduled-runner refresh of currentAveragePrice / LevePrice* / LeveProfit* columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# Row 15
$ws.Range("H15").Value = 1319.5
$ws.Range("I15").Value = 1319.5
$ws.Range("K15").Value = 3958.5
$ws.Range("M15").Value = -3789.5

# Row 18
$ws.Range("H18").Value = 845.9
$ws.Range("I18").Value = 308
$ws.Range("K18").Value = 308
$ws.Range("M18").Value = -24

# Row 40
$ws.Range("H40").Value = 4923.75
$ws.Range("I40").Value = 3071.7144
$ws.Range("J40").Value = 7516.6
$ws.Range("K40").Value = 3071.7144
$ws.Range("L40").Value = 7516.6
$ws.Range("M40").Value = -2896.7144
$ws.Range("N40").Value = -7866.6

# Row 62
$ws.Range("H62").Value = 6582663.5
$ws.Range("I62").Value = 15627025
$ws.Range("J62").Value = 4945.636
$ws.Range("K62").Value = 15627025
$ws.Range("L62").Value = 4945.636
$ws.Range("M62").Value = -15626401
$ws.Range("N62").Value = -6193.636

# Row 65
$ws.Range("H65").Value = 6582663.5
$ws.Range("I65").Value = 15627025
$ws.Range("J65").Value = 4945.636
$ws.Range("K65").Value = 78135125
$ws.Range("L65").Value = 24728.18
$ws.Range("M65").Value = -78132005
$ws.Range("N65").Value = -30968.18

# Row 76
$ws.Range("H76").Value = 62506300
$ws.Range("I76").Value = 5399.6
$ws.Range("J76").Value = 90915800
$ws.Range("K76").Value = 5399.6
$ws.Range("L76").Value = 90915800
$ws.Range("M76").Value = -5084.6
$ws.Range("N76").Value = -90916430

# Row 79
$ws.Range("H79").Value = 62506300
$ws.Range("I79").Value = 5399.6
$ws.Range("J79").Value = 90915800
$ws.Range("K79").Value = 5399.6
$ws.Range("L79").Value = 90915800
$ws.Range("M79").Value = -4307.6
$ws.Range("N79").Value = -90917984

# Row 132
$ws.Range("H132").Value = 4127.659
$ws.Range("I132").Value = 4287.579
$ws.Range("K132").Value = 12862.737
$ws.Range("M132").Value = -10332.737

# Row 141
$ws.Range("H141").Value = 2222.72
$ws.Range("I141").Value = 2181
$ws.Range("K141").Value = 6543
$ws.Range("M141").Value = -1363

$ws = $wb.Worksheets("ARM")
# Row 63
$ws.Range("H63").Value = 8364.637000000001
$ws.Range("I63").Value = 7001.8335
$ws.Range("K63").Value = 7001.8335
$ws.Range("M63").Value = -6315.8335

# Row 66
$ws.Range("H66").Value = 8364.637000000001
$ws.Range("I66").Value = 7001.8335
$ws.Range("K66").Value = 35009.1675
$ws.Range("M66").Value = -31577.1675

# Row 74
$ws.Range("H74").Value = 258782.69
$ws.Range("I74").Value = 456478.9
$ws.Range("J74").Value = 2940.5293
$ws.Range("K74").Value = 456478.9
$ws.Range("L74").Value = 2940.5293
$ws.Range("M74").Value = -455604.9
$ws.Range("N74").Value = -4688.5293

# Row 77
$ws.Range("H77").Value = 258782.69
$ws.Range("I77").Value = 456478.9
$ws.Range("J77").Value = 2940.5293
$ws.Range("K77").Value = 2282394.5
$ws.Range("L77").Value = 14702.6465
$ws.Range("M77").Value = -2278026.5
$ws.Range("N77").Value = -23438.6465

# Row 82
$ws.Range("H82").Value = 31198.5
$ws.Range("J82").Value = 31598
$ws.Range("L82").Value = 31598
$ws.Range("N82").Value = -32320

# Row 85
$ws.Range("H85").Value = 31198.5
$ws.Range("J85").Value = 31598
$ws.Range("L85").Value = 31598
$ws.Range("N85").Value = -34094

# Row 132
$ws.Range("H132").Value = 5792.7334
$ws.Range("I132").Value = 4704.25
$ws.Range("J132").Value = 7969.7
$ws.Range("K132").Value = 14112.75
$ws.Range("L132").Value = 23909.1
$ws.Range("M132").Value = -11582.75
$ws.Range("N132").Value = -28969.1

$ws = $wb.Worksheets("BSM")
# Row 74
$ws.Range("H74").Value = 28559.8
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 28559.8
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 28559.8
$ws.Range("N74").Value = -30431.8

# Row 77
$ws.Range("H77").Value = 28559.8
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 28559.8
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 85679.39999999999
$ws.Range("N77").Value = -95039.39999999999

# Row 94
$ws.Range("H94").Value = 2583.3333
$ws.Range("I94").Value = 3375
$ws.Range("K94").Value = 3375
$ws.Range("M94").Value = -2924

# Row 135
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws = $wb.Worksheets("CRP")
# Row 41
$ws.Range("H41").Value = 31936.176
$ws.Range("I41").Value = 10687.5
$ws.Range("J41").Value = 50823.89
$ws.Range("K41").Value = 10687.5
$ws.Range("L41").Value = 50823.89
$ws.Range("M41").Value = -10259.5
$ws.Range("N41").Value = -51679.89

# Row 86
$ws.Range("H86").Value = 6947.154
$ws.Range("I86").Value = 2764.75
$ws.Range("K86").Value = 2764.75
$ws.Range("M86").Value = -1641.75

# Row 89
$ws.Range("H89").Value = 6947.154
$ws.Range("I89").Value = 2764.75
$ws.Range("K89").Value = 13823.75
$ws.Range("M89").Value = -8207.75

# Row 103
$ws.Range("H103").Value = 21424.75
$ws.Range("I103").Value = 21424.75
$ws.Range("K103").Value = 21424.75
$ws.Range("M103").Value = -20252.75

# Row 107
$ws.Range("H107").Value = 1166.2858
$ws.Range("I107").Value = 690.53845
$ws.Range("K107").Value = 690.53845
$ws.Range("M107").Value = 1229.46155

# Row 122
$ws.Range("H122").Value = 3971.4285
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550

$ws = $wb.Worksheets("CUL")
# Row 80
$ws.Range("H80").Value = 4419.8
$ws.Range("J80").Value = 4419.8
$ws.Range("L80").Value = 13259.4
$ws.Range("N80").Value = -15131.4

# Row 83
$ws.Range("H83").Value = 4419.8
$ws.Range("J83").Value = 4419.8
$ws.Range("L83").Value = 39778.2
$ws.Range("N83").Value = -49138.2

# Row 128
$ws.Range("H128").Value = 148979.5
$ws.Range("I128").Value = 148979.5
$ws.Range("K128").Value = 446938.5
$ws.Range("M128").Value = -441958.5

$ws = $wb.Worksheets("GSM")
# Row 80
$ws.Range("H80").Value = 971864.3
$ws.Range("J80").Value = 1267969.5
$ws.Range("L80").Value = 1267969.5
$ws.Range("N80").Value = -1269965.5

# Row 83
$ws.Range("H83").Value = 971864.3
$ws.Range("J83").Value = 1267969.5
$ws.Range("L83").Value = 6339847.5
$ws.Range("N83").Value = -6349831.5

# Row 97
$ws.Range("H97").Value = 900.36
$ws.Range("I97").Value = 944.0952
$ws.Range("K97").Value = 944.0952
$ws.Range("M97").Value = -448.0952

# Row 122
$ws.Range("H122").Value = 561118.25
$ws.Range("I122").Value = 658279.0600000001
$ws.Range("J122").Value = 10540.333
$ws.Range("K122").Value = 1974837.18
$ws.Range("L122").Value = 31620.999
$ws.Range("M122").Value = -1972387.18
$ws.Range("N122").Value = -36520.999

# Row 132
$ws.Range("H132").Value = 41585.23
$ws.Range("I132").Value = 3012.9524
$ws.Range("J132").Value = 203588.8
$ws.Range("K132").Value = 9038.8572
$ws.Range("L132").Value = 610766.3999999999
$ws.Range("M132").Value = -6508.8572
$ws.Range("N132").Value = -615826.3999999999

$ws = $wb.Worksheets("LTW")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# Row 46
$ws.Range("H46").Value = 3400.1765
$ws.Range("I46").Value = 3000.5789
$ws.Range("J46").Value = 3906.3333
$ws.Range("K46").Value = 3000.5789
$ws.Range("L46").Value = 3906.3333
$ws.Range("M46").Value = -2812.5789
$ws.Range("N46").Value = -4282.3333

# Row 93
$ws.Range("H93").Value = 2614.0667
$ws.Range("I93").Value = 2007.2858
$ws.Range("K93").Value = 2007.2858
$ws.Range("M93").Value = -759.2858000000001

# Row 132
$ws.Range("H132").Value = 7024.6
$ws.Range("I132").Value = 4582.6665
$ws.Range("J132").Value = 8071.143
$ws.Range("K132").Value = 13747.9995
$ws.Range("L132").Value = 24213.429
$ws.Range("M132").Value = -11217.9995
$ws.Range("N132").Value = -29273.429

$ws = $wb.Worksheets("WVR")
# Row 136
$ws.Range("H136").Value = 638430.5
$ws.Range("I136").Value = 1671800.9
$ws.Range("K136").Value = 5015402.699999999
$ws.Range("M136").Value = -5012852.699999999
